# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-price refresh to the Garuda_Profits workbook.
# For each affected leve row, updates the price/profit columns (H:N) with the
# latest values; one row (ARM!N97) loses its HQ-profit figure, and one row
# (BSM!N94) gains one, matching the source data refresh exactly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 16483.883
$ws.Range("I15").Value = 16483.883
$ws.Range("K15").Value = 49451.649
$ws.Range("M15").Value = -49282.649

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 640
$ws.Range("I49").Value = 440
$ws.Range("J49").Value = 760
$ws.Range("K49").Value = 1320
$ws.Range("L49").Value = 2280
$ws.Range("M49").Value = -1184
$ws.Range("N49").Value = -2552

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1891.0646
$ws.Range("I116").Value = 1286.0714
$ws.Range("J116").Value = 2389.2942
$ws.Range("K116").Value = 1286.0714
$ws.Range("L116").Value = 2389.2942
$ws.Range("M116").Value = 2155.9286
$ws.Range("N116").Value = -9273.2942

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 90910050
$ws.Range("I125").Value = 200000480
$ws.Range("J125").Value = 1349.3334
$ws.Range("K125").Value = 1800004320
$ws.Range("L125").Value = 12144.0006
$ws.Range("M125").Value = -1800001860
$ws.Range("N125").Value = -17064.0006

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1262.3143
$ws.Range("I137").Value = 989.3
$ws.Range("J137").Value = 2900.4
$ws.Range("K137").Value = 2967.9
$ws.Range("L137").Value = 8701.200000000001
$ws.Range("M137").Value = -417.8999999999996
$ws.Range("N137").Value = -13801.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 742.8570999999999
$ws.Range("I74").Value = 405
$ws.Range("J74").Value = 1418.5714
$ws.Range("K74").Value = 405
$ws.Range("L74").Value = 1418.5714
$ws.Range("M74").Value = 469
$ws.Range("N74").Value = -3166.5714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 742.8570999999999
$ws.Range("I77").Value = 405
$ws.Range("J77").Value = 1418.5714
$ws.Range("K77").Value = 2025
$ws.Range("L77").Value = 7092.857
$ws.Range("M77").Value = 2343
$ws.Range("N77").Value = -15828.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1420.125
$ws.Range("I97").Value = 1420.125
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1420.125
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -924.125
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 165.4
$ws.Range("I80").Value = 346.8
$ws.Range("J80").Value = 74.7
$ws.Range("K80").Value = 346.8
$ws.Range("L80").Value = 74.7
$ws.Range("M80").Value = 651.2
$ws.Range("N80").Value = -2070.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 165.4
$ws.Range("I83").Value = 346.8
$ws.Range("J83").Value = 74.7
$ws.Range("K83").Value = 1734
$ws.Range("L83").Value = 373.5
$ws.Range("M83").Value = 3258
$ws.Range("N83").Value = -10357.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1860.2307
$ws.Range("I86").Value = 1687
$ws.Range("K86").Value = 1687
$ws.Range("M86").Value = -564

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1860.2307
$ws.Range("I89").Value = 1687
$ws.Range("K89").Value = 8435
$ws.Range("M89").Value = -2819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 713.75
$ws.Range("I94").Value = 694.375
$ws.Range("J94").Value = 752.5
$ws.Range("K94").Value = 694.375
$ws.Range("L94").Value = 752.5
$ws.Range("M94").Value = -243.375
$ws.Range("N94").Value = -1654.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2158.6177
$ws.Range("I105").Value = 1613.6842
$ws.Range("J105").Value = 2848.8667
$ws.Range("K105").Value = 1613.6842
$ws.Range("L105").Value = 2848.8667
$ws.Range("M105").Value = 133.3158000000001
$ws.Range("N105").Value = -6342.8667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 681.875
$ws.Range("I105").Value = 556.1539
$ws.Range("J105").Value = 830.4545000000001
$ws.Range("K105").Value = 556.1539
$ws.Range("L105").Value = 830.4545000000001
$ws.Range("M105").Value = 1190.8461
$ws.Range("N105").Value = -4324.4545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1867409.6
$ws.Range("I132").Value = 1262.3489
$ws.Range("J132").Value = 5210923.5
$ws.Range("K132").Value = 3787.0467
$ws.Range("L132").Value = 15632770.5
$ws.Range("M132").Value = -1257.0467
$ws.Range("N132").Value = -15637830.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1138.5
$ws.Range("I134").Value = 1101.5
$ws.Range("J134").Value = 1249.5
$ws.Range("K134").Value = 3304.5
$ws.Range("L134").Value = 3748.5
$ws.Range("M134").Value = -769.5
$ws.Range("N134").Value = -8818.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2433.3333
$ws.Range("I81").Value = 300
$ws.Range("J81").Value = 3500
$ws.Range("K81").Value = 900
$ws.Range("L81").Value = 10500
$ws.Range("M81").Value = 223
$ws.Range("N81").Value = -12746

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2433.3333
$ws.Range("I84").Value = 300
$ws.Range("J84").Value = 3500
$ws.Range("K84").Value = 2700
$ws.Range("L84").Value = 31500
$ws.Range("M84").Value = 2916
$ws.Range("N84").Value = -42732

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 770.99
$ws.Range("I131").Value = 371.25
$ws.Range("J131").Value = 805.75
$ws.Range("K131").Value = 1113.75
$ws.Range("L131").Value = 2417.25
$ws.Range("M131").Value = 3926.25
$ws.Range("N131").Value = -12497.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 32131548
$ws.Range("I70").Value = 44415490
$ws.Range("J70").Value = 4323.6924
$ws.Range("K70").Value = 44415490
$ws.Range("L70").Value = 4323.6924
$ws.Range("M70").Value = -44415220
$ws.Range("N70").Value = -4863.6924

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 32131548
$ws.Range("I73").Value = 44415490
$ws.Range("J73").Value = 4323.6924
$ws.Range("K73").Value = 44415490
$ws.Range("L73").Value = 4323.6924
$ws.Range("M73").Value = -44414554
$ws.Range("N73").Value = -6195.6924

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3973.0557
$ws.Range("I80").Value = 2938.75
$ws.Range("J80").Value = 4800.5
$ws.Range("K80").Value = 2938.75
$ws.Range("L80").Value = 4800.5
$ws.Range("M80").Value = -1940.75
$ws.Range("N80").Value = -6796.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3973.0557
$ws.Range("I83").Value = 2938.75
$ws.Range("J83").Value = 4800.5
$ws.Range("K83").Value = 14693.75
$ws.Range("L83").Value = 24002.5
$ws.Range("M83").Value = -9701.75
$ws.Range("N83").Value = -33986.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1074.3889
$ws.Range("I97").Value = 979.9167
$ws.Range("J97").Value = 1263.3334
$ws.Range("K97").Value = 979.9167
$ws.Range("L97").Value = 1263.3334
$ws.Range("M97").Value = -483.9167
$ws.Range("N97").Value = -2255.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 244.66667
$ws.Range("I55").Value = 259
$ws.Range("J55").Value = 233.2
$ws.Range("K55").Value = 259
$ws.Range("L55").Value = 233.2
$ws.Range("M55").Value = -86
$ws.Range("N55").Value = -579.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5079.98
$ws.Range("I132").Value = 5881.6055
$ws.Range("K132").Value = 17644.8165
$ws.Range("M132").Value = -15114.8165
